$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# "Volume 32   Number  12" -> "...13"
$ws.Range("A8").Value = "Volume 32   Number  13"
# "Report Covering the Week  3/17/2025  Through  3/23/2025" -> new dates
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# --- Row 15 (Murder) ---
# G15/H15 go from numbers (1 / -100) to the blank-marker text cells "0" / "***.*"
$ws.Range("C15").Copy($ws.Range("G15"))
$ws.Range("E15").Copy($ws.Range("H15"))
# L15 goes from text "***.*" to a real percentage number
$ws.Range("L15").Value = 100
$ws.Range("L15").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("M15").Value = 0

# --- Row 16 (Rape) ---
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = -61.111111111111
$ws.Range("L16").Value = -63.157894736842
$ws.Range("M16").Value = -69.565217391304
$ws.Range("N16").Value = -92.473118279569

# --- Row 17 (Robbery) ---
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 28
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 40
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = 86.666666666666
$ws.Range("N17").Value = 3.703703703703

# --- Row 18 (Fel. Assault) ---
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 80
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 22.222222222222
$ws.Range("I18").Value = 69
$ws.Range("J18").Value = 59
$ws.Range("K18").Value = 16.949152542372
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = 4.545454545454
$ws.Range("N18").Value = -71.369294605809

# --- Row 19 (Burglary) ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 102
$ws.Range("J19").Value = 112
$ws.Range("K19").Value = -8.928571428571
$ws.Range("L19").Value = -44.262295081967
$ws.Range("M19").Value = 10.869565217391
$ws.Range("N19").Value = -19.685039370078

# --- Row 20 (Gr. Larceny) ---
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -5.263157894736
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = -11.111111111111
$ws.Range("L20").Value = 47.368421052631
$ws.Range("M20").Value = 80.645161290322
$ws.Range("N20").Value = -92.679738562091

# --- Row 21 (G.L.A. / TOTAL) ---
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 11.764705882352
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = 8.433734939759
$ws.Range("I21").Value = 265
$ws.Range("J21").Value = 275
$ws.Range("K21").Value = -3.636363636363
$ws.Range("L21").Value = -25.352112676056
$ws.Range("M21").Value = 15.720524017467
$ws.Range("N21").Value = -78.901273885350

# --- Row 24 (Petit Larceny) ---
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 8.333333333333
$ws.Range("F24").Value = 46
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = 21.052631578947
$ws.Range("I24").Value = 134
$ws.Range("J24").Value = 144
$ws.Range("K24").Value = -6.944444444444
$ws.Range("L24").Value = 1.515151515151
$ws.Range("M24").Value = 6.349206349206

# --- Row 25 (Retail Theft) ---
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 20
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = -33.333333333333
$ws.Range("L25").Value = -20

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 20
$ws.Range("H26").Value = 42.857142857142
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 44
$ws.Range("K26").Value = -18.181818181818
$ws.Range("L26").Value = -33.333333333333
$ws.Range("M26").Value = 9.090909090909

# --- Row 27 (UCR Rape*) ---
$ws.Range("C15").Copy($ws.Range("G27"))
$ws.Range("E15").Copy($ws.Range("H27"))
$ws.Range("L27").Value = 50

# --- Row 29 (Shooting Vic.) ---
$ws.Range("C15").Copy($ws.Range("C29"))

# --- Row 30 (Shooting Inc.) ---
$ws.Range("C15").Copy($ws.Range("C30"))

# --- Row 31 (Hate Crimes) ---
$ws.Range("F31").Value = 1

Write-Output "done"
